$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.8647
$ws.Range("C2").Value = 0.5053
$ws.Range("D2").Value = 0.70103
$ws.Range("E2").Value = 1.0061
$ws.Range("F2").Value = 1.03005
$ws.Range("G2").Value = 0.6788999999999999
$ws.Range("H2").Value = 0.8031700000000001

$ws.Range("B3").Value = 0.7764
$ws.Range("C3").Value = 0.5053
$ws.Range("D3").Value = 0.8340300000000002
$ws.Range("E3").Value = 1.0144
$ws.Range("F3").Value = 0.9470500000000001
$ws.Range("G3").Value = 0.7672
$ws.Range("H3").Value = 0.8031700000000001

$ws.Range("B4").Value = 0.7847
$ws.Range("C4").Value = 0.5053
$ws.Range("D4").Value = 0.8340300000000002
$ws.Range("E4").Value = 1.0061
$ws.Range("F4").Value = 0.9770500000000001
$ws.Range("G4").Value = 0.6871999999999999
$ws.Range("H4").Value = 0.8031700000000001
